$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 10
$ws.Range("E4").Value = 9
$ws.Range("E9").Value = 10
$ws.Range("E15").Value = 10
